$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values that look numeric (e.g. "27.00", "65.02")
# but must remain plain text, exactly as authored (preserving trailing zeros,
# thousand-dot grouping, etc.). Force the whole column to Text format before
# writing any values so Excel doesn't silently coerce them into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '64.070.08'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '2.762.69'
$ws.Range("E3").Value = '  -0.46%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '576.81'
$ws.Range("E5").Value = '  -1.36%  '
$ws.Range("D6").Value = '159.57'
$ws.Range("E6").Value = '  -1.00%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '0.602'
$ws.Range("E8").Value = '  -3.24%  '
$ws.Range("E9").Value = '  -3.30%  '
$ws.Range("E10").Value = '  +3.44%  '
$ws.Range("E11").Value = '  -14.63%  '
$ws.Range("D12").Value = '0.386'
$ws.Range("E12").Value = '  -3.21%  '
$ws.Range("D13").Value = '3.251.11'
$ws.Range("E13").Value = '  -0.64%  '
$ws.Range("D14").Value = '27.00'
$ws.Range("E14").Value = '  -1.66%  '
$ws.Range("D15").Value = '63.663.14'
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("E16").Value = '  -5.04%  '
$ws.Range("D17").Value = '2.768.21'
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("D18").Value = '12.12'
$ws.Range("E18").Value = '  -1.94%  '
$ws.Range("D19").Value = '4.85'
$ws.Range("E19").Value = '  -3.19%  '
$ws.Range("D20").Value = '361.18'
$ws.Range("E20").Value = '  -1.74%  '
$ws.Range("D21").Value = '6.67'
$ws.Range("E21").Value = '  -5.45%  '
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("D23").Value = '0.528'
$ws.Range("E23").Value = '  -6.35%  '
$ws.Range("D24").Value = '65.02'
$ws.Range("E24").Value = '  -3.56%  '
$ws.Range("E25").Value = '  -4.33%  '
$ws.Range("D26").Value = '8.54'
$ws.Range("E26").Value = '  -2.04%  '
$ws.Range("E27").Value = '  +0.22%  '
$ws.Range("D28").Value = '0.0₃0907'
$ws.Range("E28").Value = '  -6.10%  '
$ws.Range("E29").Value = '  -0.73%  '
$ws.Range("D30").Value = '1.96'
$ws.Range("E30").Value = '  -3.63%  '
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("D32").Value = '170.23'
$ws.Range("E32").Value = '  -1.33%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").Value = '4.93'
$ws.Range("E33").Value = '  -3.00%  '
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").Value = '20.18'
$ws.Range("E34").Value = '  -3.25%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '1.47'
$ws.Range("E35").Value = '  -0.36%  '
$ws.Range("B36").Value = 'USDe'
$ws.Range("C36").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D36").Value = '0.998'
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("E37").Value = '  -1.47%  '
$ws.Range("E38").Value = '  -2.71%  '
$ws.Range("D39").Value = '349.04'
$ws.Range("E39").Value = '  +2.00%  '
$ws.Range("D40").Value = '6.29'
$ws.Range("E40").Value = '  +0.62%  '
$ws.Range("E41").Value = '  -2.60%  '
$ws.Range("D42").Value = '39.09'
$ws.Range("E42").Value = '  -2.07%  '
$ws.Range("E43").Value = '  -4.52%  '
$ws.Range("D44").Value = '21.79'
$ws.Range("E44").Value = '  -4.35%  '
$ws.Range("D45").Value = '0.0587'
$ws.Range("E45").Value = '  -3.70%  '
$ws.Range("D46").Value = '138.06'
$ws.Range("E46").Value = '  -0.58%  '
$ws.Range("D47").Value = '0.632'
$ws.Range("E47").Value = '  -3.52%  '
$ws.Range("E48").Value = '  -2.84%  '
$ws.Range("E49").Value = '  -1.71%  '
$ws.Range("D50").Value = '0.998'
$ws.Range("E50").Value = '  -0.03%  '
$ws.Range("D51").Value = '11.06'
$ws.Range("E51").Value = '  +0.25%  '
